# Updated cryptos list on Wed Nov 27 14:50:39 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto table, and reorders a handful of coin rows whose relative ranking
# changed (EthereumClassic/PolygonEcosystemToken, ARBITRUM/Kaspa,
# VeChain/WhiteBITCoin).
#
# Column D values are written with a leading apostrophe so Excel keeps
# them as text (matching the original inlineStr cell type) instead of
# silently re-interpreting strings such as "1.00" or "94.760.19" as
# numbers and dropping significant digits / separators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = "'94.760.19"
$ws.Range('E2').Value = '  +1.48%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = "'3.513.61"
$ws.Range('E3').Value = '  +4.64%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  +0.13%  '

# Row 5 - Solana
$ws.Range('D5').Value = "'238.71"
$ws.Range('E5').Value = '  +2.86%  '

# Row 6 - BNB
$ws.Range('D6').Value = "'629.62"
$ws.Range('E6').Value = '  +1.42%  '

# Row 7 - XRP
$ws.Range('E7').Value = '  +4.75%  '

# Row 8 - Dogecoin
$ws.Range('D8').Value = "'0.398"
$ws.Range('E8').Value = '  +3.26%  '

# Row 9 - USDC
$ws.Range('E9').Value = '  +0.09%  '

# Row 10 - Cardano
$ws.Range('E10').Value = '  +7.87%  '

# Row 11 - LidoStakedEther
$ws.Range('D11').Value = "'3.510.62"
$ws.Range('E11').Value = '  +4.60%  '

# Row 12 - Avalanche
$ws.Range('D12').Value = "'43.56"
$ws.Range('E12').Value = '  +2.99%  '

# Row 13 - TRON
$ws.Range('E13').Value = '  +4.80%  '

# Row 14 - Toncoin
$ws.Range('D14').Value = "'6.29"
$ws.Range('E14').Value = '  +4.75%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range('D15').Value = "'4.174.16"
$ws.Range('E15').Value = '  +5.04%  '

# Row 16 - WrappedBTC
$ws.Range('D16').Value = "'94.546.52"
$ws.Range('E16').Value = '  +1.53%  '

# Row 17 - ShibaInu
$ws.Range('D17').Value = "'0.0000252"
$ws.Range('E17').Value = '  +3.64%  '

# Row 18 - Polkadot
$ws.Range('D18').Value = "'8.36"
$ws.Range('E18').Value = '  +4.45%  '

# Row 19 - WrappedEther
$ws.Range('D19').Value = "'3.509.18"

# Row 20 - Uniswap
$ws.Range('D20').Value = "'12.82"
$ws.Range('E20').Value = '  +14.45%  '

# Row 21 - Chainlink
$ws.Range('D21').Value = "'18.02"
$ws.Range('E21').Value = '  +3.49%  '

# Row 22 - Stellar
$ws.Range('D22').Value = "'0.500"
$ws.Range('E22').Value = '  +9.25%  '

# Row 23 - BitcoinCash
$ws.Range('D23').Value = "'519.21"
$ws.Range('E23').Value = '  +5.71%  '

# Row 24 - SuiNetwork
$ws.Range('E24').Value = '  +1.53%  '

# Row 25 - NEARProtocol
$ws.Range('E25').Value = '  +10.05%  '

# Row 26 - PEPE
$ws.Range('D26').Value = "'0.0000187"
$ws.Range('E26').Value = '  +2.79%  '

# Row 27 - Litecoin
$ws.Range('D27').Value = "'96.45"
$ws.Range('E27').Value = '  +7.55%  '

# Row 28 - Aptos
$ws.Range('D28').Value = "'12.27"
$ws.Range('E28').Value = '  +5.39%  '

# Row 29 - PancakeSwap
$ws.Range('D29').Value = "'2.94"
$ws.Range('E29').Value = '  +10.55%  '

# Row 30 - InternetComputer(DFINITY)
$ws.Range('D30').Value = "'11.55"
$ws.Range('E30').Value = '  +3.10%  '

# Row 31 - Hedera
$ws.Range('E31').Value = '  +3.59%  '

# Row 32 - Dai
$ws.Range('D32').Value = "'1.00"
$ws.Range('E32').Value = '  +0.04%  '

# Row 33 - Cronos
$ws.Range('E33').Value = '  +5.26%  '

# Row 34 - Binance-PegBSC-USD
$ws.Range('D34').Value = "'0.995"
$ws.Range('E34').Value = '  -0.22%  '

# Row 35 - now PolygonEcosystemToken (was EthereumClassic)
$ws.Range('B35').Value = 'PolygonEcosystemToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D35').Value = "'0.563"
$ws.Range('E35').Value = '  +6.11%  '

# Row 36 - now EthereumClassic (was PolygonEcosystemToken)
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = "'29.82"
$ws.Range('E36').Value = '  +4.69%  '

# Row 37 - Bittensor
$ws.Range('D37').Value = "'583.66"
$ws.Range('E37').Value = '  +10.74%  '

# Row 38 - Fetch.AI
$ws.Range('E38').Value = '  +6.44%  '

# Row 39 - RenderToken
$ws.Range('E39').Value = '  +2.95%  '

# Row 40 - USDe
$ws.Range('E40').Value = '  +0.03%  '

# Row 41 - now Kaspa (was ARBITRUM)
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = "'0.150"
$ws.Range('E41').Value = '  +2.22%  '

# Row 42 - now ARBITRUM (was Kaspa)
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = "'0.921"
$ws.Range('E42').Value = '  +4.27%  '

# Row 43 - now WhiteBITCoin (was VeChain)
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').Value = "'23.76"
$ws.Range('E43').Value = '  -1.12%  '

# Row 44 - now VeChain (was WhiteBITCoin)
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = "'0.0424"
$ws.Range('E44').Value = '  +4.95%  '

# Row 45 - ImmutableX
$ws.Range('E45').Value = '  +1.48%  '

# Row 46 - Filecoin
$ws.Range('D46').Value = "'5.56"

# Row 47 - MantraDAO
$ws.Range('E47').Value = '  -0.31%  '

# Row 48 - Stacks
$ws.Range('E48').Value = '  +2.03%  '

# Row 49 - OKB
$ws.Range('D49').Value = "'53.94"
$ws.Range('E49').Value = '  +2.82%  '

# Row 50 - Cosmos
$ws.Range('E50').Value = '  +3.76%  '

# Row 51 - dogwifhat
$ws.Range('E51').Value = '  +0.80%  '
